$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: F1 changes meaning, G1/H1 are new ---
# F1 used to hold "TCV_range"; it now holds "Årsag" (reason)
$ws.Range("F1").Value = "Årsag"
$ws.Range("G1").Value = "Ny leverandør"
$ws.Range("H1").Value = "TCV_range"

# Copy the bold/bordered header formatting from an existing header cell (E1)
# onto the two newly introduced header cells so they match the rest of row 1.
$ws.Range("E1").Copy() | Out-Null
$ws.Range("G1:H1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Data rows: move the old "TCV_range" value (100000-120000) from F to H,
#     fill F with the new "Årsag" reason, and fill G ("Ny leverandør") only
#     where a new supplier was noted ---
$reasons = @{
    2  = "Utilfredshed (Service - uddyb i bemærkninger)"
    3  = "Utilfredshed (Service - uddyb i bemærkninger)"
    4  = "Anden årsag (angiv hvilken i bemærkninger)"
    5  = "Ikke oplyst"
    6  = "Ikke oplyst"
    7  = "Ikke oplyst"
    8  = "Ikke oplyst"
    9  = "Ikke oplyst"
    10 = "Ikke oplyst"
    11 = "Virksomheden lukker"
}

$newSuppliers = @{
    9 = "DataLøn"
}

for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 8).Value = "100000-120000"
    $ws.Cells.Item($r, 6).Value = $reasons[$r]
    if ($newSuppliers.ContainsKey($r)) {
        $ws.Cells.Item($r, 7).Value = $newSuppliers[$r]
    }
}
